$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply header row style (s="3") to A1:G1 ---
$ws.Range("A1:G1").Style = $ws.Range("B1").Style

# --- New data rows 7-11 ---
$ws.Cells.Item(7, 1).Value = 73
$ws.Cells.Item(7, 2).Value = "Baseline"
$ws.Cells.Item(7, 3).Value = 67.992369999999994
$ws.Cells.Item(7, 4).Value = 0.773567
$ws.Cells.Item(7, 5).Value = [double]"6.2268500000000002E-05"
$ws.Cells.Item(7, 6).Value = [double]"1.61633E-04"
$ws.Cells.Item(7, 7).Value = [double]"2.0894500000000001E-04"

$ws.Cells.Item(8, 1).Value = 73
$ws.Cells.Item(8, 2).Value = "Stressor"
$ws.Cells.Item(8, 3).Value = 70.476900000000001
$ws.Cells.Item(8, 4).Value = 0.422593
$ws.Cells.Item(8, 5).Value = [double]"2.7780600000000002E-07"
$ws.Cells.Item(8, 6).Value = [double]"1.0250999999999999E-06"
$ws.Cells.Item(8, 7).Value = [double]"2.42575E-06"

$ws.Cells.Item(9, 1).Value = 73
$ws.Cells.Item(9, 2).Value = "R1"
$ws.Cells.Item(9, 3).Value = 68.067750000000004
$ws.Cells.Item(9, 4).Value = 1.0098199999999999
$ws.Cells.Item(9, 5).Value = [double]"7.7207300000000002E-05"
$ws.Cells.Item(9, 6).Value = [double]"1.7866400000000001E-04"
$ws.Cells.Item(9, 7).Value = [double]"1.7692599999999999E-04"

$ws.Cells.Item(10, 1).Value = 73
$ws.Cells.Item(10, 2).Value = "R2"
$ws.Cells.Item(10, 3).Value = 70.51576
$ws.Cells.Item(10, 4).Value = 1.2548699999999999
$ws.Cells.Item(10, 5).Value = [double]"4.3641199999999998E-05"
$ws.Cells.Item(10, 6).Value = [double]"2.01883E-04"
$ws.Cells.Item(10, 7).Value = [double]"1.6087999999999999E-04"

$ws.Cells.Item(11, 1).Value = 73
$ws.Cells.Item(11, 2).Value = "R3"
$ws.Cells.Item(11, 3).Value = 69.109520000000003
$ws.Cells.Item(11, 4).Value = 1.8723099999999999
$ws.Cells.Item(11, 5).Value = [double]"3.6699599999999997E-05"
$ws.Cells.Item(11, 6).Value = [double]"3.4311299999999998E-04"
$ws.Cells.Item(11, 7).Value = [double]"1.8325000000000001E-04"

# --- Copy formatting/styles from row 2 (same pattern as existing data rows) ---
for ($r = 7; $r -le 11; $r++) {
    $ws.Range("B$r").Style = $ws.Range("B2").Style
    $ws.Range("C$r").Style = $ws.Range("C2").Style
    $ws.Range("D$r").Style = $ws.Range("D2").Style
    $ws.Range("F$r").Style = $ws.Range("F2").Style
    $ws.Range("G$r").Style = $ws.Range("G2").Style
}
# Column E styles vary per-row (style 1 vs style 2) to mirror the source rows
$ws.Range("E7").Style = $ws.Range("G2").Style
$ws.Range("E8").Style = $ws.Range("G2").Style
$ws.Range("E9").Style = $ws.Range("G2").Style
$ws.Range("E10").Style = $ws.Range("G2").Style
$ws.Range("E11").Style = $ws.Range("G2").Style
# Columns F/G on row 8 use style 2 (small-exponent numbers) instead of style 1
$ws.Range("F8").Style = $ws.Range("G2").Style
$ws.Range("G8").Style = $ws.Range("G2").Style

# --- Selection / view changes ---
$ws.Range("B12").Select()

$excel.ActiveWindow.Left = 368
$excel.ActiveWindow.Top = 368
$excel.ActiveWindow.Width = 14399
$excel.ActiveWindow.Height = 8272
